# Applies the GDP/Elasticities/Intensities workbook edits described in the
# commit diff: updated GDP reference-year values & growth-rate assumptions,
# removed the outdated B1 cell comment on the GDP sheet, removed an unused
# spreadsheet column on the GDP sheet, refreshed the elasticity and
# intensity uncertainty-table figures, and restored the expected sheet
# selections.

$wb = $excel.ActiveWorkbook

$wsGDP   = $wb.Worksheets.Item("GDP")
$wsElas  = $wb.Worksheets.Item("Elasticities")
$wsInten = $wb.Worksheets.Item("Intensities")

# ---------------------------------------------------------------------
# GDP sheet
# ---------------------------------------------------------------------

# The old sourcing note on B1 (GDP value derived from CRC exchange rate)
# no longer applies now that GDP is entered directly - drop the comment.
$wsGDP.Range("B1").Comment.Delete()

# 2018-2020 GDP figures are now hard values instead of derived formulas.
$wsGDP.Range("B2").Value2 = 68004
$wsGDP.Range("B3").Value2 = 70634
$wsGDP.Range("C3").ClearContents()
$wsGDP.Range("B4").Value2 = 69561
$wsGDP.Range("C4").ClearContents()

# From 2021 onward, GDP keeps growing off the prior year at a flat 3.5%
# per year (replacing the old variable growth-rate assumptions).
$wsGDP.Range("C5:C34").Value2 = 3.5
$wsGDP.Range("B5").Formula = "=B4*(1+C5/100)"
for ($r = 6; $r -le 34; $r++) {
    $prev = $r - 1
    $wsGDP.Range("B$r").Formula = "=B$prev*(1+C$r/100)"
}

# The now-unused helper column (empty, only currency-formatted) is removed.
$wsGDP.Columns("E").Delete()

# ---------------------------------------------------------------------
# Elasticities sheet - refreshed e_Passenger / e_Freight assumptions
# ---------------------------------------------------------------------

$wsElas.Range("B2").Value2 = 2.4
$wsElas.Range("C2").Value2 = 1.8

$wsElas.Range("B3").Value2 = -6.2
$wsElas.Range("C3").Value2 = -3.6

$wsElas.Range("B4").Value2 = 1.6
$wsElas.Range("C4").Value2 = 1.1

$wsElas.Range("B5").Value2 = 1.6
$wsElas.Range("C5").Value2 = 1.1

$wsElas.Range("B6").Value2 = 1.6
$wsElas.Range("C6").Value2 = 1.1

$wsElas.Range("B7").Value2 = 1.6
$wsElas.Range("C7").Value2 = 1.1

$wsElas.Range("C33").Value2 = 1

# ---------------------------------------------------------------------
# Intensities sheet - refreshed i_NT_elec / i_NT_fossil assumptions
# ---------------------------------------------------------------------

$wsInten.Range("B2").Value2 = 569.69220119404747
$wsInten.Range("C2").Value2 = 617.84772289277112

$wsInten.Range("B3").Value2 = 572.09535039782543
$wsInten.Range("C3").Value2 = 726.66152336268647

$wsInten.Range("B4").Value2 = 565.7144974914105
$wsInten.Range("C4").Value2 = 757.1376149652823

$wsInten.Range("B5").Value2 = 560.69234983613035
$wsInten.Range("C5").Value2 = 749.24155834931821

$wsInten.Range("B6").Value2 = 560.69234983613035
$wsInten.Range("C6").Value2 = 749.24155834931821

$wsInten.Range("B7").Value2 = 560.69234983613035
$wsInten.Range("C7").Value2 = 749.24155834931821

$wsInten.Range("B8").Value2 = 560.69234983613035
$wsInten.Range("C8").Value2 = 749.24155834931821

$wsInten.Range("B34").Value2 = 392.49492137322432
$wsInten.Range("C34").Value2 = 497.92920331860444

# ---------------------------------------------------------------------
# Restore the expected per-sheet selections / active tab.
# Elasticities must be selected last so it ends up as the active tab.
# ---------------------------------------------------------------------

[void]$wsGDP.Range("B1").Select()
[void]$wsInten.Range("B8:C8").Select()
[void]$wsElas.Range("B1").Select()
